$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The edit adds a new price-record row in the middle of the data
# (pushing the existing rows 799..879 down to 800..880), and appends
# one more row at the very end (881) which duplicates the data that
# is now in row 880 (it was originally row 879).
# -----------------------------------------------------------------

# 1) Insert a new row just above the old row 799, shifting everything
#    below it down by one.
$ws.Rows.Item(799).Insert()

# 2) Populate the newly inserted row 799 with its data.
$ws.Cells.Item(799, 1).Value  = 10
$ws.Cells.Item(799, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(799, 3).Value  = "La Araucanía"
$ws.Cells.Item(799, 4).Value  = 45132
$ws.Cells.Item(799, 5).Value  = 9
$ws.Cells.Item(799, 6).Value  = 100112032
$ws.Cells.Item(799, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(799, 8).Value  = "Bola 8"
$ws.Cells.Item(799, 9).Value  = "Primera"
$ws.Cells.Item(799, 10).Value = 110
$ws.Cells.Item(799, 11).Value = 20000
$ws.Cells.Item(799, 12).Value = 20000
$ws.Cells.Item(799, 13).Value = 20000
$ws.Cells.Item(799, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(799, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(799, 16).Value = 400
$ws.Cells.Item(799, 17).Value = 50
$ws.Cells.Item(799, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of
# column D (the row-insert already carries this down, but set it
# explicitly to be safe).
$ws.Cells.Item(799, 4).NumberFormat = $ws.Cells.Item(800, 4).NumberFormat

# 3) Append a new row 881 duplicating the row that is now at 880
#    (originally row 879).
$ws.Cells.Item(881, 1).Value  = 10
$ws.Cells.Item(881, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(881, 3).Value  = "La Araucanía"
$ws.Cells.Item(881, 4).Value  = 45072
$ws.Cells.Item(881, 5).Value  = 9
$ws.Cells.Item(881, 6).Value  = 100112032
$ws.Cells.Item(881, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(881, 8).Value  = "Sin especificar"
$ws.Cells.Item(881, 9).Value  = "Primera"
$ws.Cells.Item(881, 10).Value = 250
$ws.Cells.Item(881, 11).Value = 13000
$ws.Cells.Item(881, 12).Value = 13000
$ws.Cells.Item(881, 13).Value = 13000
$ws.Cells.Item(881, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(881, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(881, 16).Value = 260
$ws.Cells.Item(881, 17).Value = 50
$ws.Cells.Item(881, 18).Value = "Hortaliza"

$ws.Cells.Item(881, 4).NumberFormat = $ws.Cells.Item(880, 4).NumberFormat
